$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "31.211.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.65%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.998.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +6.74%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7842"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +66.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "257.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.85%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9985"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3540"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +23.52%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "29.47"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +34.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.38"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.17%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07051"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.56%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8722"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +19.70%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08205"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.002.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "101.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.601"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "15.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +18.36%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "274.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "31.223.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.961"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +11.90%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000007957"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.262.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9989"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.13%  "

$ws.Range("E24").Value = "  -0.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.158"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +12.99%  "

$ws.Range("E26").Value = "  +11.91%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1492"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +54.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.61%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.383"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +25.71%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.608"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.626"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.353"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.40%  "

$ws.Range("E34").Value = "  +7.13%  "

$ws.Range("E35").Value = "  +8.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.232"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7799"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +13.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.798"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.34%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02014"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.910"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.36%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.760"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "79.87"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.94%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4736"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +12.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.153"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "106.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.62%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.8491"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9986"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.756"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.970"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4349"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.88"
$ws.Range("D51").Style = "Normal"
